# Insert a new weekly record for "Puerro" (Vega Modelo de Temuco) as row 121,
# pushing the existing rows 121-151 down to 122-152.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 121; Excel shifts rows 121:151 down to 122:152
# and copies row formatting (incl. the date style on column D) from the row above.
$ws.Rows("121:121").Insert()

# Populate the new row 121 with the new daily record.
$ws.Cells.Item(121, 1).Value2  = 10                               # A Mercado ID
$ws.Cells.Item(121, 2).Value2  = "Vega Modelo de Temuco"          # B Mercado
$ws.Cells.Item(121, 3).Value2  = "La Araucanía"                   # C Región
$ws.Cells.Item(121, 4).Value2  = 44508                            # D Fecha
$ws.Cells.Item(121, 5).Value2  = 9                                # E Codreg
$ws.Cells.Item(121, 6).Value2  = 100112005                        # F Categoría ID
$ws.Cells.Item(121, 7).Value2  = "Puerro"                         # G Categoría
$ws.Cells.Item(121, 8).Value2  = "Azul de Maquehue"               # H Variedad
$ws.Cells.Item(121, 9).Value2  = "Primera"                        # I Calidad
$ws.Cells.Item(121, 10).Value2 = 80                               # J Volumen
$ws.Cells.Item(121, 11).Value2 = 7000                             # K Precio mínimo
$ws.Cells.Item(121, 12).Value2 = 7000                             # L Precio máximo
$ws.Cells.Item(121, 13).Value2 = 7000                             # M Precio promedio ponderado
$ws.Cells.Item(121, 14).Value2 = "$/docena de paquetes"           # N Unidad de comercialización
$ws.Cells.Item(121, 15).Value2 = "Provincia de Cautín"            # O Origen
$ws.Cells.Item(121, 16).Value2 = 583                              # P Precio $/Kg
$ws.Cells.Item(121, 17).Value2 = 12                               # Q Kg o Unidades
$ws.Cells.Item(121, 18).Value2 = "Hortaliza"                      # R Clasificación
